# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values, rows 2-20 (data rows for A=0..18)
$newK = @(4, 3, 7, 5, 0, 6, 6, 5, 5, 3, 8, 8, 9, 5, 2, 5, 2, 5, 4)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
